$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update the Correspond Handoff/Handback datetimes for the
# e6052173-... row (row 3) to reflect the newly generated handback report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-24 07:09:40"
$wsZhCn.Range("G3").Value = "2016-02-24 07:10:33"

# "de-de" sheet: same update for its e6052173-... row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-24 07:09:52"
$wsDeDe.Range("G3").Value = "2016-02-24 07:10:57"
